$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("M20").ClearContents() | Out-Null
$ws.Range("H35").Value = 0
$ws.Range("I35").Value = 0
$ws.Range("K35").Value = 0
$ws.Range("M35").ClearContents() | Out-Null
$ws.Range("H100").Value = 2265.5789
$ws.Range("I100").Value = 1682
$ws.Range("J100").Value = 2914
$ws.Range("K100").Value = 1682
$ws.Range("L100").Value = 2914
$ws.Range("M100").Value = -1141
$ws.Range("N100").Value = -3996
$ws.Range("H107").Value = 491.92856
$ws.Range("I107").Value = 397
$ws.Range("J107").Value = 840
$ws.Range("K107").Value = 397
$ws.Range("L107").Value = 840
$ws.Range("M107").Value = 1523
$ws.Range("N107").Value = -4680
$ws.Range("H127").Value = 940.63635
$ws.Range("I127").Value = 940.63635
$ws.Range("J127").Value = 0
$ws.Range("K127").Value = 2821.90905
$ws.Range("L127").Value = 0
$ws.Range("M127").Value = 2138.09095
$ws.Range("N127").ClearContents() | Out-Null
$ws.Range("H138").Value = 2058.1875
$ws.Range("I138").Value = 871.2778
$ws.Range("J138").Value = 3584.2144
$ws.Range("K138").Value = 2613.8334
$ws.Range("L138").Value = 10752.6432
$ws.Range("M138").Value = 2526.1666
$ws.Range("N138").Value = -21032.6432

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 630.5599999999999
$ws.Range("I32").Value = 619.1340300000001
$ws.Range("J32").Value = 1000
$ws.Range("K32").Value = 619.1340300000001
$ws.Range("L32").Value = 1000
$ws.Range("M32").Value = -332.1340300000001
$ws.Range("N32").Value = -1574
$ws.Range("H88").Value = 2736.8096
$ws.Range("I88").Value = 2398.8333
$ws.Range("J88").Value = 3187.4443
$ws.Range("K88").Value = 2398.8333
$ws.Range("L88").Value = 3187.4443
$ws.Range("M88").Value = -1992.8333
$ws.Range("N88").Value = -3999.4443
$ws.Range("H91").Value = 2736.8096
$ws.Range("I91").Value = 2398.8333
$ws.Range("J91").Value = 3187.4443
$ws.Range("K91").Value = 2398.8333
$ws.Range("L91").Value = 3187.4443
$ws.Range("M91").Value = -994.8332999999998
$ws.Range("N91").Value = -5995.4443

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1956.9584
$ws.Range("I86").Value = 1422.5294
$ws.Range("J86").Value = 3254.8572
$ws.Range("K86").Value = 1422.5294
$ws.Range("L86").Value = 3254.8572
$ws.Range("M86").Value = -299.5293999999999
$ws.Range("N86").Value = -5500.8572
$ws.Range("H89").Value = 1956.9584
$ws.Range("I89").Value = 1422.5294
$ws.Range("J89").Value = 3254.8572
$ws.Range("K89").Value = 7112.646999999999
$ws.Range("L89").Value = 16274.286
$ws.Range("M89").Value = -1496.646999999999
$ws.Range("N89").Value = -27506.286
$ws.Range("H134").Value = 892.02563
$ws.Range("I134").Value = 821.69446
$ws.Range("K134").Value = 2465.08338
$ws.Range("M134").Value = 69.91661999999997

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 33051.457
$ws.Range("I31").Value = 1077.2368
$ws.Range("J31").Value = 73552.13
$ws.Range("K31").Value = 1077.2368
$ws.Range("L31").Value = 73552.13
$ws.Range("M31").Value = -782.2367999999999
$ws.Range("N31").Value = -74142.13
$ws.Range("H34").Value = 33051.457
$ws.Range("I34").Value = 1077.2368
$ws.Range("J34").Value = 73552.13
$ws.Range("K34").Value = 1077.2368
$ws.Range("L34").Value = 73552.13
$ws.Range("M34").Value = -875.2367999999999
$ws.Range("N34").Value = -73956.13
$ws.Range("H132").Value = 16396495
$ws.Range("I132").Value = 19611116
$ws.Range("J132").Value = 1927.4
$ws.Range("K132").Value = 58833348
$ws.Range("L132").Value = 5782.200000000001
$ws.Range("M132").Value = -58830818
$ws.Range("N132").Value = -10842.2

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 476690.16
$ws.Range("I5").Value = 280.25
$ws.Range("J5").Value = 2001201.8
$ws.Range("K5").Value = 840.75
$ws.Range("L5").Value = 6003605.4
$ws.Range("M5").Value = -728.75
$ws.Range("N5").Value = -6003829.4
$ws.Range("H122").Value = 872.7406999999999
$ws.Range("I122").Value = 405.33334
$ws.Range("J122").Value = 1006.2857
$ws.Range("K122").Value = 3648.00006
$ws.Range("L122").Value = 9056.5713
$ws.Range("M122").Value = -1198.00006
$ws.Range("N122").Value = -13956.5713
$ws.Range("H135").Value = 476690.16
$ws.Range("I135").Value = 280.25
$ws.Range("J135").Value = 2001201.8
$ws.Range("K135").Value = 2522.25
$ws.Range("L135").Value = 18010816.2
$ws.Range("M135").Value = 12.75
$ws.Range("N135").Value = -18015886.2

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2756.2666
$ws.Range("I80").Value = 2090
$ws.Range("J80").Value = 3517.7144
$ws.Range("K80").Value = 2090
$ws.Range("L80").Value = 3517.7144
$ws.Range("M80").Value = -1092
$ws.Range("N80").Value = -5513.7144
$ws.Range("H83").Value = 2756.2666
$ws.Range("I83").Value = 2090
$ws.Range("J83").Value = 3517.7144
$ws.Range("K83").Value = 10450
$ws.Range("L83").Value = 17588.572
$ws.Range("M83").Value = -5458
$ws.Range("N83").Value = -27572.572
$ws.Range("H132").Value = 3053.5806
$ws.Range("I132").Value = 3451.5
$ws.Range("J132").Value = 1689.2858
$ws.Range("K132").Value = 10354.5
$ws.Range("L132").Value = 5067.857400000001
$ws.Range("M132").Value = -7824.5
$ws.Range("N132").Value = -10127.8574

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 936.8889
$ws.Range("I82").Value = 850.6667
$ws.Range("J82").Value = 980
$ws.Range("K82").Value = 850.6667
$ws.Range("L82").Value = 980
$ws.Range("M82").Value = -489.6667
$ws.Range("N82").Value = -1702
$ws.Range("H85").Value = 936.8889
$ws.Range("I85").Value = 850.6667
$ws.Range("J85").Value = 980
$ws.Range("K85").Value = 850.6667
$ws.Range("L85").Value = 980
$ws.Range("M85").Value = 397.3333
$ws.Range("N85").Value = -3476

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 2625
$ws.Range("I62").Value = 2600
$ws.Range("J62").Value = 2660
$ws.Range("K62").Value = 2600
$ws.Range("L62").Value = 2660
$ws.Range("M62").Value = -1976
$ws.Range("N62").Value = -3908
$ws.Range("H65").Value = 2625
$ws.Range("I65").Value = 2600
$ws.Range("J65").Value = 2660
$ws.Range("K65").Value = 13000
$ws.Range("L65").Value = 13300
$ws.Range("M65").Value = -9880
$ws.Range("N65").Value = -19540
$ws.Range("H81").Value = 806.75
$ws.Range("I81").Value = 771.9048
$ws.Range("J81").Value = 911.2857
$ws.Range("K81").Value = 1543.8096
$ws.Range("L81").Value = 1822.5714
$ws.Range("M81").Value = -482.8096
$ws.Range("N81").Value = -3944.5714
$ws.Range("H84").Value = 806.75
$ws.Range("I84").Value = 771.9048
$ws.Range("J84").Value = 911.2857
$ws.Range("K84").Value = 7719.048000000001
$ws.Range("L84").Value = 9112.857
$ws.Range("M84").Value = -2415.048000000001
$ws.Range("N84").Value = -19720.857
$ws.Range("H132").Value = 2474.1365
$ws.Range("I132").Value = 3536.9211
$ws.Range("J132").Value = 1031.7858
$ws.Range("K132").Value = 10610.7633
$ws.Range("L132").Value = 3095.3574
$ws.Range("M132").Value = -8080.763300000001
$ws.Range("N132").Value = -8155.357400000001
